$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Target layout (row 1 = headers, row 2 = one data record). Every populated
# cell is formatted as Text ("@") before the value is written so that
# numeric-looking values (e.g. "17", "1987", "12345678") are stored as
# shared strings rather than numbers - matching the source workbook.
# ---------------------------------------------------------------------------

$cells = [ordered]@{
    "A1" = "Key"
    "B1" = "FirstName"
    "C1" = "LastName"
    "D1" = "MartialStatus"
    "E1" = "Hobbies"
    "F1" = "Country"
    "G1" = "BirthMonth"
    "H1" = "BirthDay"
    "I1" = "BirthYear"
    "J1" = "Phone"
    "K1" = "UserName"
    "L1" = "Email"
    "M1" = "Picture"
    "N1" = "Description"
    "O1" = "Password"
    "P1" = "ConfirmPassword"

    "A2" = "RegisterWithoutLastName"
    "B2" = "Iliya"
    "D2" = "1,2,3"
    "E2" = "1,2,3"
    "F2" = "Bulgaria"
    "G2" = "3"
    "H2" = "17"
    "I2" = "1987"
    "J2" = "0897675645"
    "K2" = "lichkata456"
    "L2" = "lichkata456@abv.bg"
    "M2" = "C:\Users\Iliya\Desktop\photo.jpeg"
    "N2" = "ALA BALA"
    "O2" = "12345678"
    "P2" = "12345678"
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
}

# Email cell becomes a live mailto hyperlink (adopts the built-in Hyperlink
# style automatically on top of the Text format already applied above).
$ws.Hyperlinks.Add($ws.Range("L2"), "mailto:lichkata456@abv.bg") | Out-Null

# Column sizing to fit the new content (values chosen so the engine's
# char-width -> xlsx "width" conversion lands as close as possible on the
# authored widths).
$ws.Columns("D:D").ColumnWidth = 11.833333333333334
$ws.Columns("K:K").ColumnWidth = 10.0
$ws.Columns("L:L").ColumnWidth = 18.0
$ws.Columns("M:M").ColumnWidth = 31.666666666666668
$ws.Columns("N:N").ColumnWidth = 10.333333333333334
$ws.Columns("O:O").ColumnWidth = 8.666666666666666
$ws.Columns("P:P").ColumnWidth = 15.833333333333334

# Final selection / view mirrors the saved state in the workbook.
$ws.Range("P2").Select() | Out-Null

Write-Host "UserData sheet populated"
